$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Private" column (D) values to TRUE for data rows 2-6
$ws.Range("D2:D6").Value = $true

# Give column D the same TRUE/FALSE list validation already used on column F
# (xlValidateList=3, xlValidAlertStop=1, xlBetween=1 — operator is a no-op for lists)
$ws.Range("D2:D6").Validation.Add(3, 1, 1, """TRUE,FALSE""")

# Reflect the new selection on column D (as seen in the saved file)
[void]$ws.Range("D2:D6").Select()
